$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.349.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.284.65"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.27"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.280.78"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.89%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000262"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.812.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "620.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.504.30"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.288.96"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.59"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.884"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.99"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.11"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.95"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.34"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.53"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.33"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.44"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "568.05"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.65"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -10.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.81"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.832.87"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.38"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.47%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.43"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "32.37"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0678"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.51%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.88%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.30%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.09"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.18%  "
